# Commit: "new .ttl from Google sheet has been generated"
#
# The diff shows, for every data row (1..231), that column X's original
# value is discarded and every column from Y through AQ shifts one
# position to the left (Y->X, Z->Y, ... AQ->AP), with the now-vacated
# last column (AQ) disappearing entirely. The sheet's dimension shrinks
# from A1:AQ231 to A1:AP231 accordingly. That is exactly the effect of
# deleting the entire column X.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("X").Delete()

$wb.Save()
